$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Статистика по городам")

# Widen column A (15.22 -> 16.22) to match new layout.
# Note: the host quantizes ColumnWidth to whole pixels internally, so the
# exact value 16.22 isn't reachable; 15.33 lands on the closest achievable
# stored width (16.1667 char-units, the nearest representable value to 16.22).
$ws.Columns.Item(1).ColumnWidth = 15.33

# --- New rows 12-24: stamp cell formatting (style) by copying from existing template cells ---
# A/B columns for rows 12-24 use the same style as existing data rows (s=2)
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A12:B12").PasteSpecial(-4122) | Out-Null
$ws.Range("A13:B13").PasteSpecial(-4122) | Out-Null
$ws.Range("A14:B14").PasteSpecial(-4122) | Out-Null
$ws.Range("A15:B15").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:B16").PasteSpecial(-4122) | Out-Null
$ws.Range("A17:B17").PasteSpecial(-4122) | Out-Null
$ws.Range("A18:B18").PasteSpecial(-4122) | Out-Null
$ws.Range("A19:B19").PasteSpecial(-4122) | Out-Null
$ws.Range("A20:B20").PasteSpecial(-4122) | Out-Null
$ws.Range("A21:B21").PasteSpecial(-4122) | Out-Null
$ws.Range("A22:B22").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:B23").PasteSpecial(-4122) | Out-Null
$ws.Range("A24:B24").PasteSpecial(-4122) | Out-Null

# D column for rows 12-17 uses style s=2 (same template as D2)
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null

# E column for rows 12-17 uses style s=2 (NOT the percent style s=3 used by rows 2-11);
# use A2 (style s=2) as the format template instead of E2 (style s=3)
$ws.Range("A2").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Write values for all data rows ---
$ws.Range("A2").Value = "Воронеж"
$ws.Range("B2").Value = 73587
$ws.Range("D2").Value = "Москва"
$ws.Range("E2").Value = 0.2461

$ws.Range("A3").Value = "Ростов-на-Дону"
$ws.Range("B3").Value = 76579
$ws.Range("D3").Value = "Санкт-Петербург"
$ws.Range("E3").Value = 0.1086

$ws.Range("A4").Value = "Казань"
$ws.Range("B4").Value = 59991
$ws.Range("D4").Value = "Нижний Новгород"
$ws.Range("E4").Value = 0.022

$ws.Range("A5").Value = "Нижний Новгород"
$ws.Range("B5").Value = 76491
$ws.Range("D5").Value = "Ростов-на-Дону"
$ws.Range("E5").Value = 0.02

$ws.Range("A6").Value = "Новосибирск"
$ws.Range("B6").Value = 46820
$ws.Range("D6").Value = "Воронеж"
$ws.Range("E6").Value = 0.0167

$ws.Range("A7").Value = "Екатеринбург"
$ws.Range("B7").Value = 81550
$ws.Range("D7").Value = "Екатеринбург"
$ws.Range("E7").Value = 0.0305

$ws.Range("A8").Value = "Санкт-Петербург"
$ws.Range("B8").Value = 52774
$ws.Range("D8").Value = "Казань"
$ws.Range("E8").Value = 0.0269

$ws.Range("A9").Value = "Москва"
$ws.Range("B9").Value = 51234
$ws.Range("D9").Value = "Новосибирск"
$ws.Range("E9").Value = 0.0308

$ws.Range("A10").Value = "Россия"
$ws.Range("B10").Value = 48221
$ws.Range("D10").Value = "Россия"
$ws.Range("E10").Value = 0.0224

$ws.Range("A11").Value = "Самара"
$ws.Range("B11").Value = 78886
$ws.Range("D11").Value = "Самара"
$ws.Range("E11").Value = 0.0162

$ws.Range("A12").Value = "Ярославль"
$ws.Range("B12").Value = 46273
$ws.Range("D12").Value = "Краснодар"
$ws.Range("E12").Value = 0.0243

$ws.Range("A13").Value = "Краснодар"
$ws.Range("B13").Value = 82585
$ws.Range("D13").Value = "Ярославль"
$ws.Range("E13").Value = 0.0127

$ws.Range("A14").Value = "Пермь"
$ws.Range("B14").Value = 81440
$ws.Range("D14").Value = "Красноярск"
$ws.Range("E14").Value = 0.0145

$ws.Range("A15").Value = "Красноярск"
$ws.Range("B15").Value = 69883
$ws.Range("D15").Value = "Пермь"
$ws.Range("E15").Value = 0.0165

$ws.Range("A16").Value = "Волгоград"
$ws.Range("B16").Value = 46518
$ws.Range("D16").Value = "Уфа"
$ws.Range("E16").Value = 0.0152

$ws.Range("A17").Value = "Уфа"
$ws.Range("B17").Value = 79632
$ws.Range("D17").Value = "Челябинск"
$ws.Range("E17").Value = 0.0167

$ws.Range("A18").Value = "Саратов"
$ws.Range("B18").Value = 50935

$ws.Range("A19").Value = "Тула"
$ws.Range("B19").Value = 44614

$ws.Range("A20").Value = "Ижевск"
$ws.Range("B20").Value = 44443

$ws.Range("A21").Value = "Челябинск"
$ws.Range("B21").Value = 70440

$ws.Range("A22").Value = "Омск"
$ws.Range("B22").Value = 49788

$ws.Range("A23").Value = "Рязань"
$ws.Range("B23").Value = 37293

$ws.Range("A24").Value = "Тюмень"
$ws.Range("B24").Value = 78655
